$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 1117
$wsExhibit.Range("F3").Value = 245
$wsExhibit.Range("F4").Value = 243
$wsExhibit.Range("F5").Value = 1782
$wsExhibit.Range("F6").Value = 667
$wsExhibit.Range("F7").Value = 321
$wsExhibit.Range("F8").Value = 496
$wsExhibit.Range("F9").Value = 4571
$wsExhibit.Range("F10").Value = 55
$wsExhibit.Range("F13").Value = 986
$wsExhibit.Range("F14").Value = 1286
$wsExhibit.Range("F18").Value = 1815
$wsExhibit.Range("F20").Value = 44
$wsExhibit.Range("F22").Value = 25
$wsExhibit.Range("F23").Value = 655
$wsExhibit.Range("F25").Value = 302
$wsExhibit.Range("F26").Value = 32
$wsExhibit.Range("F27").Value = 2379
$wsExhibit.Range("F28").Value = 1004
$wsExhibit.Range("F29").Value = 2433
$wsExhibit.Range("F30").Value = 249
$wsExhibit.Range("F31").Value = 1121
$wsExhibit.Range("F32").Value = 594
$wsExhibit.Range("F35").Value = 423
$wsExhibit.Range("F36").Value = 1116
$wsExhibit.Range("F37").Value = 915
$wsExhibit.Range("F38").Value = 1178
$wsExhibit.Range("F39").Value = 13
$wsExhibit.Range("F40").Value = 851
$wsExhibit.Range("F42").Value = 360
$wsExhibit.Range("F43").Value = 283
$wsExhibit.Range("F44").Value = 3480

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F10").Value = 882
$wsShow.Range("F23").Value = 27

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 1117
$wsAll.Range("F3").Value = 245
$wsAll.Range("F4").Value = 243
$wsAll.Range("F6").Value = 1782
$wsAll.Range("F7").Value = 667
$wsAll.Range("F8").Value = 321
$wsAll.Range("F9").Value = 496
$wsAll.Range("F10").Value = 4571
$wsAll.Range("F11").Value = 55
$wsAll.Range("F15").Value = 1286
$wsAll.Range("F18").Value = 1815
$wsAll.Range("F20").Value = 44
$wsAll.Range("F23").Value = 882
$wsAll.Range("F25").Value = 26
$wsAll.Range("F27").Value = 302
$wsAll.Range("F28").Value = 2379
$wsAll.Range("F31").Value = 1004
$wsAll.Range("F33").Value = 2433
$wsAll.Range("F34").Value = 1121
$wsAll.Range("F35").Value = 594
$wsAll.Range("F37").Value = 1116
$wsAll.Range("F38").Value = 915
$wsAll.Range("F40").Value = 1179
$wsAll.Range("F41").Value = 851
$wsAll.Range("F44").Value = 360
$wsAll.Range("F46").Value = 27
$wsAll.Range("F47").Value = 283
$wsAll.Range("F48").Value = 3480
